$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cx3cl1"
$ws.Cells.Item(2, 3).Value = "Itgb3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 4.550961
$ws.Cells.Item(2, 8).Value = 13.652883
$ws.Cells.Item(2, 9).Value = 0.3980483771262702
$ws.Cells.Item(2, 10).Value = 0.3980483771262702
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.214110666666667
$ws.Cells.Item(2, 14).Value = 21.642332
$ws.Cells.Item(2, 15).Value = 0.4688823795981188
$ws.Cells.Item(2, 16).Value = 0.4688823795981188
$ws.Cells.Item(2, 17).Value = 32.831136293684
$ws.Cells.Item(2, 18).Value = 295.480226643156
$ws.Cells.Item(2, 19).Value = 0.186637870262135
$ws.Cells.Item(2, 20).Value = 0.186637870262135

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cx3cl1"
$ws.Cells.Item(3, 3).Value = "Itgb3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 4.550961
$ws.Cells.Item(3, 8).Value = 13.652883
$ws.Cells.Item(3, 9).Value = 0.3980483771262702
$ws.Cells.Item(3, 10).Value = 0.3980483771262702
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.110350666666666
$ws.Cells.Item(3, 14).Value = 21.331052
$ws.Cells.Item(3, 15).Value = 0.4621384803214003
$ws.Cells.Item(3, 16).Value = 0.4621384803214003
$ws.Cells.Item(3, 17).Value = 32.358928580324
$ws.Cells.Item(3, 18).Value = 291.230357222916
$ws.Cells.Item(3, 19).Value = 0.1839534720995341
$ws.Cells.Item(3, 20).Value = 0.1839534720995341

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cx3cl1"
$ws.Cells.Item(4, 3).Value = "Itgb3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 4.550961
$ws.Cells.Item(4, 8).Value = 13.652883
$ws.Cells.Item(4, 9).Value = 0.3980483771262702
$ws.Cells.Item(4, 10).Value = 0.3980483771262702
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.061296333333333
$ws.Cells.Item(4, 14).Value = 3.183889
$ws.Cells.Item(4, 15).Value = 0.06897914008048092
$ws.Cells.Item(4, 16).Value = 0.06897914008048092
$ws.Cells.Item(4, 17).Value = 4.829918222442999
$ws.Cells.Item(4, 18).Value = 43.46926400198699
$ws.Cells.Item(4, 19).Value = 0.02745703476460109
$ws.Cells.Item(4, 20).Value = 0.02745703476460109

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cx3cl1"
$ws.Cells.Item(5, 3).Value = "Itgb3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.449245
$ws.Cells.Item(5, 8).Value = 16.347735
$ws.Cells.Item(5, 9).Value = 0.4766165055717775
$ws.Cells.Item(5, 10).Value = 0.4766165055717775
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 7.214110666666667
$ws.Cells.Item(5, 14).Value = 21.642332
$ws.Cells.Item(5, 15).Value = 0.4688823795981188
$ws.Cells.Item(5, 16).Value = 0.4688823795981188
$ws.Cells.Item(5, 17).Value = 39.31145647978
$ws.Cells.Item(5, 18).Value = 353.80310831802
$ws.Cells.Item(5, 19).Value = 0.2234770812882351
$ws.Cells.Item(5, 20).Value = 0.2234770812882351

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Cx3cl1"
$ws.Cells.Item(6, 3).Value = "Itgb3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.449245
$ws.Cells.Item(6, 8).Value = 16.347735
$ws.Cells.Item(6, 9).Value = 0.4766165055717775
$ws.Cells.Item(6, 10).Value = 0.4766165055717775
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.110350666666666
$ws.Cells.Item(6, 14).Value = 21.331052
$ws.Cells.Item(6, 15).Value = 0.4621384803214003
$ws.Cells.Item(6, 16).Value = 0.4621384803214003
$ws.Cells.Item(6, 17).Value = 38.74604281858
$ws.Cells.Item(6, 18).Value = 348.71438536722
$ws.Cells.Item(6, 19).Value = 0.2202628275810375
$ws.Cells.Item(6, 20).Value = 0.2202628275810375

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Cx3cl1"
$ws.Cells.Item(7, 3).Value = "Itgb3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.449245
$ws.Cells.Item(7, 8).Value = 16.347735
$ws.Cells.Item(7, 9).Value = 0.4766165055717775
$ws.Cells.Item(7, 10).Value = 0.4766165055717775
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.061296333333333
$ws.Cells.Item(7, 14).Value = 3.183889
$ws.Cells.Item(7, 15).Value = 0.06897914008048092
$ws.Cells.Item(7, 16).Value = 0.06897914008048092
$ws.Cells.Item(7, 17).Value = 5.783263737934999
$ws.Cells.Item(7, 18).Value = 52.049373641415
$ws.Cells.Item(7, 19).Value = 0.03287659670250496
$ws.Cells.Item(7, 20).Value = 0.03287659670250496

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Cx3cl1"
$ws.Cells.Item(8, 3).Value = "Itgb3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.432979666666667
$ws.Cells.Item(8, 8).Value = 4.298939
$ws.Cells.Item(8, 9).Value = 0.1253351173019523
$ws.Cells.Item(8, 10).Value = 0.1253351173019523
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 7.214110666666667
$ws.Cells.Item(8, 14).Value = 21.642332
$ws.Cells.Item(8, 15).Value = 0.4688823795981188
$ws.Cells.Item(8, 16).Value = 0.4688823795981188
$ws.Cells.Item(8, 17).Value = 10.33767389841644
$ws.Cells.Item(8, 18).Value = 93.03906508574799
$ws.Cells.Item(8, 19).Value = 0.05876742804774876
$ws.Cells.Item(8, 20).Value = 0.05876742804774876

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Cx3cl1"
$ws.Cells.Item(9, 3).Value = "Itgb3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.432979666666667
$ws.Cells.Item(9, 8).Value = 4.298939
$ws.Cells.Item(9, 9).Value = 0.1253351173019523
$ws.Cells.Item(9, 10).Value = 0.1253351173019523
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 7.110350666666666
$ws.Cells.Item(9, 14).Value = 21.331052
$ws.Cells.Item(9, 15).Value = 0.4621384803214003
$ws.Cells.Item(9, 16).Value = 0.4621384803214003
$ws.Cells.Item(9, 17).Value = 10.18898792820311
$ws.Cells.Item(9, 18).Value = 91.700891353828
$ws.Cells.Item(9, 19).Value = 0.05792218064082869
$ws.Cells.Item(9, 20).Value = 0.05792218064082869

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Cx3cl1"
$ws.Cells.Item(10, 3).Value = "Itgb3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.432979666666667
$ws.Cells.Item(10, 8).Value = 4.298939
$ws.Cells.Item(10, 9).Value = 0.1253351173019523
$ws.Cells.Item(10, 10).Value = 0.1253351173019523
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.061296333333333
$ws.Cells.Item(10, 14).Value = 3.183889
$ws.Cells.Item(10, 15).Value = 0.06897914008048092
$ws.Cells.Item(10, 16).Value = 0.06897914008048092
$ws.Cells.Item(10, 17).Value = 1.520816065974555
$ws.Cells.Item(10, 18).Value = 13.687344593771
$ws.Cells.Item(10, 19).Value = 0.008645508613374877
$ws.Cells.Item(10, 20).Value = 0.008645508613374877
